$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 38.93002066666667
$ws.Range("H2").Value = 116.790062
$ws.Range("I2").Value = 0.7610372167397395
$ws.Range("J2").Value = 0.7610372167397393
$ws.Range("M2").Value = 0.5229673333333335
$ws.Range("N2").Value = 1.568902
$ws.Range("O2").Value = 0.01867311683630173
$ws.Range("P2").Value = 0.01867311683630173
$ws.Range("Q2").Value = 20.35912909465823
$ws.Range("R2").Value = 183.232161851924
$ws.Range("S2").Value = 0.01421093686495504
$ws.Range("T2").Value = 0.01421093686495503

$ws.Range("G3").Value = 38.93002066666667
$ws.Range("H3").Value = 116.790062
$ws.Range("I3").Value = 0.7610372167397395
$ws.Range("J3").Value = 0.7610372167397393
$ws.Range("O3").Value = 0.1760454316627222
$ws.Range("P3").Value = 0.1760454316627221
$ws.Range("Q3").Value = 191.940729615016
$ws.Range("R3").Value = 1727.466566535144
$ws.Range("S3").Value = 0.1339771253323441
$ws.Range("T3").Value = 0.133977125332344

$ws.Range("G4").Value = 38.93002066666667
$ws.Range("H4").Value = 116.790062
$ws.Range("I4").Value = 0.7610372167397395
$ws.Range("J4").Value = 0.7610372167397393
$ws.Range("M4").Value = 1.127819333333333
$ws.Range("N4").Value = 3.383458
$ws.Range("O4").Value = 0.04027001466294246
$ws.Range("P4").Value = 0.04027001466294246
$ws.Range("Q4").Value = 43.90602995493289
$ws.Range("R4").Value = 395.1542695943961
$ws.Range("S4").Value = 0.03064697987715423
$ws.Range("T4").Value = 0.03064697987715423

$ws.Range("G5").Value = 38.93002066666667
$ws.Range("H5").Value = 116.790062
$ws.Range("I5").Value = 0.7610372167397395
$ws.Range("J5").Value = 0.7610372167397393
$ws.Range("M5").Value = 21.42523899999999
$ws.Range("N5").Value = 64.27571699999999
$ws.Range("O5").Value = 0.7650114368380336
$ws.Range("P5").Value = 0.7650114368380336
$ws.Range("Q5").Value = 834.0849970582725
$ws.Range("R5").Value = 7506.764973524453
$ws.Range("S5").Value = 0.5822021746652861
$ws.Range("T5").Value = 0.582202174665286

$ws.Range("I6").Value = 0.1914142145281647
$ws.Range("J6").Value = 0.1914142145281647
$ws.Range("M6").Value = 0.5229673333333335
$ws.Range("N6").Value = 1.568902
$ws.Range("O6").Value = 0.01867311683630173
$ws.Range("P6").Value = 0.01867311683630173
$ws.Range("Q6").Value = 5.120678224944446
$ws.Range("R6").Value = 46.0861040245
$ws.Range("S6").Value = 0.003574299992013344
$ws.Range("T6").Value = 0.003574299992013343

$ws.Range("I7").Value = 0.1914142145281647
$ws.Range("J7").Value = 0.1914142145281647
$ws.Range("O7").Value = 0.1760454316627222
$ws.Range("P7").Value = 0.1760454316627221
$ws.Range("S7").Value = 0.03369759802299167
$ws.Range("T7").Value = 0.03369759802299165

$ws.Range("I8").Value = 0.1914142145281647
$ws.Range("J8").Value = 0.1914142145281647
$ws.Range("M8").Value = 1.127819333333333
$ws.Range("N8").Value = 3.383458
$ws.Range("O8").Value = 0.04027001466294246
$ws.Range("P8").Value = 0.04027001466294246
$ws.Range("Q8").Value = 11.04313698727778
$ws.Range("R8").Value = 99.3882328855
$ws.Range("S8").Value = 0.007708253225744808
$ws.Range("T8").Value = 0.007708253225744807

$ws.Range("I9").Value = 0.1914142145281647
$ws.Range("J9").Value = 0.1914142145281647
$ws.Range("M9").Value = 21.42523899999999
$ws.Range("N9").Value = 64.27571699999999
$ws.Range("O9").Value = 0.7650114368380336
$ws.Range("P9").Value = 0.7650114368380336
$ws.Range("Q9").Value = 209.7870131050833
$ws.Range("R9").Value = 1888.08311794575
$ws.Range("S9").Value = 0.1464340632874149
$ws.Range("T9").Value = 0.1464340632874149

$ws.Range("G10").Value = 1.794146
$ws.Range("H10").Value = 5.382438
$ws.Range("I10").Value = 0.03507349482179579
$ws.Range("J10").Value = 0.03507349482179579
$ws.Range("M10").Value = 0.5229673333333335
$ws.Range("N10").Value = 1.568902
$ws.Range("O10").Value = 0.01867311683630173
$ws.Range("P10").Value = 0.01867311683630173
$ws.Range("Q10").Value = 0.938279749230667
$ws.Range("R10").Value = 8.444517743076002
$ws.Range("S10").Value = 0.0006549314666648166
$ws.Range("T10").Value = 0.0006549314666648164

$ws.Range("G11").Value = 1.794146
$ws.Range("H11").Value = 5.382438
$ws.Range("I11").Value = 0.03507349482179579
$ws.Range("J11").Value = 0.03507349482179579
$ws.Range("O11").Value = 0.1760454316627222
$ws.Range("P11").Value = 0.1760454316627221
$ws.Range("Q11").Value = 8.845864614984002
$ws.Range("R11").Value = 79.61278153485601
$ws.Range("S11").Value = 0.006174528535823291
$ws.Range("T11").Value = 0.006174528535823289

$ws.Range("G12").Value = 1.794146
$ws.Range("H12").Value = 5.382438
$ws.Range("I12").Value = 0.03507349482179579
$ws.Range("J12").Value = 0.03507349482179579
$ws.Range("M12").Value = 1.127819333333333
$ws.Range("N12").Value = 3.383458
$ws.Range("O12").Value = 0.04027001466294246
$ws.Range("P12").Value = 0.04027001466294246
$ws.Range("Q12").Value = 2.023472545622667
$ws.Range("R12").Value = 18.211252910604
$ws.Range("S12").Value = 0.001412410150754353
$ws.Range("T12").Value = 0.001412410150754353

$ws.Range("G13").Value = 1.794146
$ws.Range("H13").Value = 5.382438
$ws.Range("I13").Value = 0.03507349482179579
$ws.Range("J13").Value = 0.03507349482179579
$ws.Range("M13").Value = 21.42523899999999
$ws.Range("N13").Value = 64.27571699999999
$ws.Range("O13").Value = 0.7650114368380336
$ws.Range("P13").Value = 0.7650114368380336
$ws.Range("Q13").Value = 38.44000685089399
$ws.Range("R13").Value = 345.960061658046
$ws.Range("S13").Value = 0.02683162466855333
$ws.Range("T13").Value = 0.02683162466855333

$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.6381486666666666
$ws.Range("H14").Value = 1.914446
$ws.Range("I14").Value = 0.01247507391030006
$ws.Range("J14").Value = 0.01247507391030006
$ws.Range("M14").Value = 0.5229673333333335
$ws.Range("N14").Value = 1.568902
$ws.Range("O14").Value = 0.01867311683630173
$ws.Range("P14").Value = 0.01867311683630173
$ws.Range("Q14").Value = 0.3337309064768889
$ws.Range("R14").Value = 3.003578158292
$ws.Range("S14").Value = 0.0002329485126685326
$ws.Range("T14").Value = 0.0002329485126685325

$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.6381486666666666
$ws.Range("H15").Value = 1.914446
$ws.Range("I15").Value = 0.01247507391030006
$ws.Range("J15").Value = 0.01247507391030006
$ws.Range("O15").Value = 0.1760454316627222
$ws.Range("P15").Value = 0.1760454316627221
$ws.Range("Q15").Value = 3.146330738728
$ws.Range("R15").Value = 28.316976648552
$ws.Range("S15").Value = 0.002196179771563138
$ws.Range("T15").Value = 0.002196179771563138

$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.6381486666666666
$ws.Range("H16").Value = 1.914446
$ws.Range("I16").Value = 0.01247507391030006
$ws.Range("J16").Value = 0.01247507391030006
$ws.Range("M16").Value = 1.127819333333333
$ws.Range("N16").Value = 3.383458
$ws.Range("O16").Value = 0.04027001466294246
$ws.Range("P16").Value = 0.04027001466294246
$ws.Range("Q16").Value = 0.7197164038075554
$ws.Range("R16").Value = 6.477447634268
$ws.Range("S16").Value = 0.0005023714092890746
$ws.Range("T16").Value = 0.0005023714092890746

$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.6381486666666666
$ws.Range("H17").Value = 1.914446
$ws.Range("I17").Value = 0.01247507391030006
$ws.Range("J17").Value = 0.01247507391030006
$ws.Range("M17").Value = 21.42523899999999
$ws.Range("N17").Value = 64.27571699999999
$ws.Range("O17").Value = 0.7650114368380336
$ws.Range("P17").Value = 0.7650114368380336
$ws.Range("Q17").Value = 13.67248770086466
$ws.Range("R17").Value = 123.052389307782
$ws.Range("S17").Value = 0.00954357421677932
$ws.Range("T17").Value = 0.00954357421677932

